$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.019.28'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.304.89'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'300.12"
$ws.Range('D5').ClearFormats()
$ws.Range('D6').Value = "'97.84"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').Value = "'35.97"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('D11').Value = "'0.0791"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = "'18.13"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.119"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '2.664.25'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '2.303.46'
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = "'0.782"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').Value = '42.932.82'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = "'12.73"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.90%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').Value = "'67.90"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').Value = "'236.11"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.46"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').Value = "'25.40"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.49%  '
$ws.Range('D29').Value = "'165.53"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = "'9.06"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('D32').Value = "'33.34"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('D36').Value = "'17.00"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.25%  '
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').Value = "'2.74"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').Value = '2.014.38'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').Value = "'0.0283"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'10.09"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = "'2.13"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = "'17.62"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').Value = "'2.96"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = "'53.96"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').Value = '2.532.04'
$ws.Range('E51').Value = '  +0.11%  '
